# Add "Yichun's test" (with spell-check proofErr markup around the
# possessive "Yichun's") to the single empty paragraph in the document.
$d = $word.ActiveDocument

$apostrophe = [char]0x2019   # U+2019 RIGHT SINGLE QUOTATION MARK ("smart quote")

$xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Yichun' + $apostrophe + 's</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> test</w:t></w:r>' +
       '</w:p>' +
       '</pkg:xmlData>'

$d.Paragraphs(1).Range.InsertXML($xml)
